$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("pH 9.2")
$ws.Rows("1:2").Delete()
